$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/volume updates (GitHub Actions scrape refresh)
$ws.Range("D2").Value = '41.782.36'
$ws.Range("E2").Value = '  +5.61%  '
$ws.Range("D3").Value = '2.222.87'
$ws.Range("E3").Value = '  +2.85%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = "'231.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").Value = "'60.62"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.13%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.401"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").Value = "'58.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("D11").Value = "'0.0887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.69%  '
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("D13").Value = '2.554.78'
$ws.Range("E13").Value = '  +2.90%  '
$ws.Range("D14").Value = "'15.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").Value = "'21.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("D17").Value = "'5.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.78%  '
$ws.Range("D18").Value = '2.222.67'
$ws.Range("E18").Value = '  +2.65%  '
$ws.Range("D19").Value = '41.692.18'
$ws.Range("E19").Value = '  +5.38%  '
$ws.Range("D20").Value = "'72.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.58%  '
$ws.Range("D21").Value = '0.0₃0890'
$ws.Range("E21").Value = '  +1.83%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = "'249.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.67%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").Value = "'9.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.45%  '
$ws.Range("E28").Value = '  +2.77%  '
$ws.Range("D29").Value = "'167.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").Value = "'19.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.38%  '
$ws.Range("D31").Value = "'1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").Value = "'0.121"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.22%  '
$ws.Range("D34").Value = "'4.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.35%  '
$ws.Range("D35").Value = "'4.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("D36").Value = "'0.0623"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("E37").Value = '  -4.60%  '
$ws.Range("E38").Value = '  -4.07%  '
$ws.Range("D39").Value = "'2.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.43%  '
$ws.Range("D40").Value = "'0.000257"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +28.91%  '
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("E42").Value = '  +5.50%  '
$ws.Range("D43").Value = "'4.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("D44").Value = "'8.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.59%  '
$ws.Range("D45").Value = "'0.0979"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.94%  '
$ws.Range("D46").Value = "'1.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("D47").Value = "'98.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.74%  '
$ws.Range("D48").Value = '1.468.03'
$ws.Range("E48").Value = '  -2.91%  '
$ws.Range("D49").Value = "'16.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.36%  '
$ws.Range("E50").Value = '  +0.14%  '

# Row 51: coin swapped out (ARBITRUM -> MultiversX)
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = "'52.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.43%  '
